$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.843082308769226
$ws.Range("B1").Value = 3.960108757019043
$ws.Range("C1").Value = 2.475994825363159
$ws.Range("D1").Value = 0.895289421081543
$ws.Range("E1").Value = 0.5871200561523438
